$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so numeric-looking
# strings (e.g. "1.015", "0.4970") are preserved exactly as text,
# matching how the source data was stored (inline strings).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.374.09'
$ws.Range("E2").Value = '  +1.70%  '
$ws.Range("D3").Value = '1.843.23'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("D4").Value = '1.015'
$ws.Range("E4").Value = '  +1.35%  '
$ws.Range("D5").Value = '315.81'
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("D6").Value = '1.013'
$ws.Range("E6").Value = '  +1.16%  '
$ws.Range("D7").Value = '0.4739'
$ws.Range("E7").Value = '  +1.31%  '
$ws.Range("D8").Value = '0.3703'
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").Value = '0.07454'
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("D10").Value = '0.8870'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("D11").Value = '20.52'
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("D12").Value = '1.835.99'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '0.07381'
$ws.Range("E13").Value = '  +4.29%  '
$ws.Range("D14").Value = '5.488'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("D15").Value = '93.51'
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").Value = '6.591'
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '1.015'
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").Value = '0.000008847'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").Value = '1.012'
$ws.Range("E19").Value = '  +1.09%  '
$ws.Range("D20").Value = '14.87'
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = '27.384.48'
$ws.Range("E21").Value = '  +1.56%  '
$ws.Range("D22").Value = '5.343'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("E23").Value = '  +1.37%  '
$ws.Range("D24").Value = '2.074.46'
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = '1.910'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").Value = '152.71'
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("D27").Value = '18.71'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '2.187'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").Value = '5.302'
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("D30").Value = '118.11'
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("D31").Value = '0.08971'
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("D32").Value = '0.7624'
$ws.Range("E32").Value = '  -1.06%  '
$ws.Range("D33").Value = '1.179'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("D34").Value = '4.568'
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("D35").Value = '2.950'
$ws.Range("E35").Value = '  +1.67%  '
$ws.Range("E36").Value = '  +1.22%  '
$ws.Range("D37").Value = '1.108'
$ws.Range("E37").Value = '  +2.12%  '
$ws.Range("D38").Value = '0.05369'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = '0.01967'
$ws.Range("E39").Value = '  -0.17%  '
$ws.Range("D40").Value = '3.003'
$ws.Range("E40").Value = '  +2.13%  '
$ws.Range("D41").Value = '7.367'
$ws.Range("D42").Value = '2.420'
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("D43").Value = '0.5366'
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("D44").Value = '0.1671'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '8.569'
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").Value = '0.4970'
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").Value = '10.58'
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("D48").Value = '1.014'
$ws.Range("E48").Value = '  +1.30%  '
$ws.Range("D49").Value = '1.686'
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("D50").Value = '104.67'
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  +0.74%  '
